# Applies the cryptocurrency price/volume updates described in the commit diff
# ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Some new values look like plain numbers (e.g. "1.00", "234.66") and Excel
    # would silently coerce them into numeric cells, stripping meaningful trailing
    # zeros. Force a text number format first so the literal string is kept, then
    # restore the default "Normal" style so no visible formatting change remains.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '43.503.13'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.361.29'
$ws.Range('E3').Value = '  +5.31%  '
$ws.Range('E4').Value = '  +0.16%  '
Set-TextValue 'D5' '234.66'
$ws.Range('E5').Value = '  +1.94%  '
Set-TextValue 'D6' '0.651'
$ws.Range('E6').Value = '  +1.27%  '
Set-TextValue 'D7' '72.41'
$ws.Range('E7').Value = '  +14.05%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +14.13%  '
$ws.Range('E10').Value = '  +3.69%  '
Set-TextValue 'D11' '27.26'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '2.715.81'
$ws.Range('E12').Value = '  +5.56%  '
$ws.Range('E13').Value = '  +2.40%  '
Set-TextValue 'D14' '16.46'
$ws.Range('E14').Value = '  +8.31%  '
$ws.Range('E15').Value = '  +4.95%  '
Set-TextValue 'D16' '0.872'
$ws.Range('E16').Value = '  +6.03%  '
$ws.Range('D17').Value = '2.357.22'
$ws.Range('E17').Value = '  +5.60%  '
$ws.Range('D18').Value = '43.422.82'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('E20').Value = '  +3.38%  '
$ws.Range('E21').Value = '  +5.27%  '
Set-TextValue 'D22' '251.38'
$ws.Range('E22').Value = '  +2.09%  '
Set-TextValue 'D23' '3.85'
$ws.Range('E23').Value = '  +3.64%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +3.39%  '
Set-TextValue 'D26' '10.12'
$ws.Range('E26').Value = '  +4.56%  '
Set-TextValue 'D27' '2.25'
$ws.Range('E27').Value = '  +1.20%  '
Set-TextValue 'D28' '22.67'
$ws.Range('E28').Value = '  +4.68%  '
Set-TextValue 'D29' '172.86'
$ws.Range('E29').Value = '  -0.41%  '
Set-TextValue 'D30' '1.55'
$ws.Range('E30').Value = '  +10.08%  '
$ws.Range('E31').Value = '  +3.78%  '
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('E33').Value = '  +2.18%  '
Set-TextValue 'D34' '0.0697'
$ws.Range('E34').Value = '  +3.42%  '
$ws.Range('E35').Value = '  +3.90%  '
Set-TextValue 'D36' '3.76'
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('E37').Value = '  +5.11%  '
Set-TextValue 'D38' '2.44'
$ws.Range('E38').Value = '  +7.92%  '
Set-TextValue 'D39' '0.0258'
$ws.Range('E39').Value = '  +3.54%  '
Set-TextValue 'D40' '19.61'
$ws.Range('E40').Value = '  +16.16%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D41' '1.00'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '8.94'
$ws.Range('E42').Value = '  +3.91%  '
Set-TextValue 'D43' '100.11'
$ws.Range('E43').Value = '  +4.01%  '
$ws.Range('E44').Value = '  +10.34%  '
Set-TextValue 'D45' '4.52'
$ws.Range('E45').Value = '  +2.21%  '
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('D48').Value = '1.446.07'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').Value = '2.587.45'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('E51').Value = '  -2.10%  '
